$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated week 5 points for each athlete
$ws.Range("B2").Value = 1587.8
$ws.Range("B3").Value = 1478.3
$ws.Range("B4").Value = 1276.9000000000001
$ws.Range("B5").Value = 1149
$ws.Range("B6").Value = 884.2
$ws.Range("B7").Value = 865.8
$ws.Range("B8").Value = 836.8
$ws.Range("B9").Value = 743.8
$ws.Range("B10").Value = 451.8
$ws.Range("B11").Value = 141

# Update the selected cell in the sheet view
$ws.Range("D8").Select()
